$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string for A3 (Branch value)
$ws.Range("A3").Value = "CIMS.CAN.MB.Residential.Dwellings.Lighting"

# Replace shared formulas in row 3 (M3:W3) with static values
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0.5
$ws.Range("Q3").Value = 0.0001
$ws.Range("R3").Value = 0.0001
$ws.Range("S3").Value = 0.0001
$ws.Range("T3").Value = 0.0001
$ws.Range("U3").Value = 0.0001
$ws.Range("V3").Value = 0.0001
$ws.Range("W3").Value = 0.0001

# Update the selected range in the sheet view
$ws.Range("A1:X4").Select()
